$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45172 -> 2023-09-03).
# Update all data rows (C2:C123) to the new date serial 45175 (2023-09-06).
$ws.Range("C2:C123").Value = 45175
